# Automatic map update (2025-10-13 07:27:20)
# Removes three resolved/duplicate case rows from the "AYKO" sheet and lets
# every row below them shift up, matching the refreshed export.
#   - Caso 6475  (Av Amancio Alcorta 3570)  -> old row 63
#   - Caso -602  (Agustin de vedia 2111)    -> old row 65
#   - Caso 7296  (VEDIA, AGUSTIN DE 2130)   -> old row 69

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so the remaining row numbers referenced below
# stay valid while each EntireRow delete shifts subsequent rows upward.
$ws.Rows.Item(69).Delete()
$ws.Rows.Item(65).Delete()
$ws.Rows.Item(63).Delete()
